# TC03_Verify_BLP_Solutions_ContactUS.xlsx
#
# Update the expected "Contact Us" header text on the Testdata sheet
# (previously "Contact Kaman Industrial Technologies"), and leave the
# workbook with the Testdata sheet active/selected (cell B12), while
# the first sheet keeps a lingering selection at C13.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("TC03_Verify_BLP_Sol_Contact_etc")
$ws2 = $wb.Worksheets.Item("Testdata")

# ContactUSHeader data row: value column changes to "Contact Us"
$ws2.Range("B11").Value = "Contact Us"

# Leave a selection on sheet 1 (it is no longer the active/selected tab)
$ws1.Activate()
$ws1.Range("C13").Select() | Out-Null

# Make Testdata the active sheet with its own selection - this is the
# tab that ends up marked tabSelected / activeTab in the saved file.
$ws2.Activate()
$ws2.Range("B12").Select() | Out-Null
